$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new numeric-looking text must stay text (avoid float coercion)
$textCells = @("D5", "D8", "D9", "D19", "D23", "D24", "D26", "D34", "D41", "D44", "D46", "D47", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "29.465.13"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "1.605.79"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "212.72"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("E6").Value = "  +7.02%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").Value = "26.79"
$ws.Range("E8").Value = "  +6.37%  "
$ws.Range("D9").Value = "43.63"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  +2.61%  "
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "1.836.50"
$ws.Range("E13").Value = "  +2.69%  "
$ws.Range("D14").Value = "1.580.74"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "29.478.70"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").Value = "240.29"
$ws.Range("E19").Value = "  +5.02%  "
$ws.Range("E20").Value = "  +3.62%  "
$ws.Range("D21").Value = "0.0₃0690"
$ws.Range("E21").Value = "  +1.58%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("D23").Value = "3.99"
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("D24").Value = "9.20"
$ws.Range("E24").Value = "  +1.86%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "154.40"
$ws.Range("E26").Value = "  +1.95%  "
$ws.Range("E27").Value = "  +4.77%  "
$ws.Range("E28").Value = "  +3.29%  "
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  +2.42%  "
$ws.Range("E32").Value = "  +1.21%  "
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "3.09"
$ws.Range("E34").Value = "  +3.59%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "1.412.80"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("E37").Value = "  +3.99%  "
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("D41").Value = "0.538"
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("E43").Value = "  +5.84%  "
$ws.Range("D44").Value = "0.798"
$ws.Range("E44").Value = "  +3.37%  "
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("D46").Value = "52.72"
$ws.Range("E46").Value = "  +21.69%  "
$ws.Range("D47").Value = "65.87"
$ws.Range("E47").Value = "  +2.80%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "1.745.36"
$ws.Range("E49").Value = "  +2.84%  "
$ws.Range("D50").Value = "0.859"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").Value = "86.63"
$ws.Range("E51").Value = "  +1.79%  "

Write-Host "Updated cryptos sheet."
